$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B (ASIN) to hold Week_Start_Date
$ws.Columns.Item(2).Insert()

# Header for the newly inserted column
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week start dates (Sundays), one per data row (rows 2..17)
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2

    # Shorten the week label, e.g. "W01" -> "W1"
    $ws.Cells.Item($row, 1).Value = "W" + ($i + 1)

    # Write the week start date as text (quote-prefixed so Excel keeps it as a string)
    $ws.Cells.Item($row, 2).Value = "'" + $weekStartDates[$i]

    # is_holiday_week (column J after the insert) should be boolean, not numeric
    $ws.Cells.Item($row, 10).Value = $false
}

Write-Output "done"
